# Daily attendance processing - 2025-12-04 06:36:58
# Normalize the "Recorded By" column (G) so that the "System" user is
# always listed first in the comma-separated list of recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($value -eq "system, System, backup@backdoor.com") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
}
